# Updates cryptocurrency price/volume data in the active worksheet,
# matching the values captured in the latest data refresh.
# Column map: B=Coin, C=Link, D=Price, E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '43.742.55'
$ws.Range("E2").Value = '  +0.62%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.330.65'
$ws.Range("E3").Value = '  +4.32%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.10%  '

# Row 5: Solana
$ws.Range("B5").Value = 'Solana'
$ws.Range("C5").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D5").Value = '''95.90'
$ws.Range("E5").Value = '  +2.15%  '

# Row 6: BNB
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '''270.65'
$ws.Range("E6").Value = '  +0.30%  '

# Row 7: XRP
$ws.Range("D7").Value = '''0.627'
$ws.Range("E7").Value = '  +0.73%  '

# Row 8: USDC
$ws.Range("E8").Value = '  -0.08%  '

# Row 9: Cardano
$ws.Range("E9").Value = '  -0.30%  '

# Row 10: Avalanche
$ws.Range("D10").Value = '''45.43'
$ws.Range("E10").Value = '  -2.38%  '

# Row 11: Dogecoin
$ws.Range("D11").Value = '''0.0946'
$ws.Range("E11").Value = '  +2.86%  '

# Row 12: Polkadot
$ws.Range("D12").Value = '''8.12'
$ws.Range("E12").Value = '  +0.38%  '

# Row 13: TRON
$ws.Range("E13").Value = '  +0.38%  '

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = '2.677.64'
$ws.Range("E14").Value = '  +4.21%  '

# Row 15: Chainlink
$ws.Range("D15").Value = '''15.64'
$ws.Range("E15").Value = '  +3.48%  '

# Row 16: Polygon
$ws.Range("D16").Value = '''0.865'
$ws.Range("E16").Value = '  +8.07%  '

# Row 17: WrappedEther
$ws.Range("D17").Value = '2.332.88'
$ws.Range("E17").Value = '  +4.31%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '43.711.44'
$ws.Range("E18").Value = '  +0.61%  '

# Row 19: ShibaInu
$ws.Range("E19").Value = '  +5.96%  '

# Row 20: Uniswap
$ws.Range("E20").Value = '  +7.05%  '

# Row 21: Litecoin
$ws.Range("D21").Value = '''72.59'
$ws.Range("E21").Value = '  +3.10%  '

# Row 22: BitcoinCash
$ws.Range("D22").Value = '''239.05'
$ws.Range("E22").Value = '  +2.80%  '

# Row 23: ImmutableX
$ws.Range("D23").Value = '''2.27'
$ws.Range("E23").Value = '  -2.33%  '

# Row 24: InternetComputer(DFINITY)
$ws.Range("D24").Value = '''9.38'
$ws.Range("E24").Value = '  +7.04%  '

# Row 25: Dai
$ws.Range("E25").Value = '  -0.07%  '

# Row 26: PancakeSwap
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").Value = '''2.53'

# Row 27: Cosmos
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '''11.43'
$ws.Range("E27").Value = '  +2.21%  '

# Row 28: WEMIXToken
$ws.Range("D28").Value = '''3.49'
$ws.Range("E28").Value = '  -2.07%  '

# Row 29: Toncoin
$ws.Range("D29").Value = '''2.30'
$ws.Range("E29").Value = '  +1.43%  '

# Row 30: EthereumClassic
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = '''22.56'
$ws.Range("E30").Value = '  +8.47%  '

# Row 31: InjectiveProtocol
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = '''38.17'
$ws.Range("E31").Value = '  -3.27%  '

# Row 32: Monero
$ws.Range("D32").Value = '''173.20'
$ws.Range("E32").Value = '  -0.02%  '

# Row 33: Hedera
$ws.Range("E33").Value = '  -2.10%  '

# Row 34: Filecoin
$ws.Range("D34").Value = '''5.49'
$ws.Range("E34").Value = '  +1.14%  '

# Row 35: Stellar
$ws.Range("D35").Value = '''0.127'
$ws.Range("E35").Value = '  +2.58%  '

# Row 36: VeChain
$ws.Range("E36").Value = '  +3.96%  '

# Row 37: Kaspa
$ws.Range("E37").Value = '  -2.21%  '

# Row 38: RenderToken
$ws.Range("D38").Value = '''4.39'
$ws.Range("E38").Value = '  +2.08%  '

# Row 39: NEARProtocol
$ws.Range("D39").Value = '''3.40'
$ws.Range("E39").Value = '  -3.65%  '

# Row 40: LidoDAOToken
$ws.Range("E40").Value = '  +10.16%  '

# Row 41: Algorand
$ws.Range("D41").Value = '''0.238'
$ws.Range("E41").Value = '  +9.65%  '

# Row 42: ARBITRUM
$ws.Range("D42").Value = '''1.39'
$ws.Range("E42").Value = '  +20.10%  '

# Row 43: Celestia
$ws.Range("D43").Value = '''12.13'

# Row 44: FraxShare
$ws.Range("D44").Value = '''9.23'
$ws.Range("E44").Value = '  +10.26%  '

# Row 45: MultiversX
$ws.Range("D45").Value = '''62.33'
$ws.Range("E45").Value = '  -0.94%  '

# Row 46: THORChain
$ws.Range("E46").Value = '  -0.01%  '

# Row 47: Cronos
$ws.Range("E47").Value = '  +4.91%  '

# Row 48: Aave
$ws.Range("D48").Value = '''100.47'
$ws.Range("E48").Value = '  +0.47%  '

# Row 49: TrustWalletToken
$ws.Range("E49").Value = '  +1.91%  '

# Row 50: RocketPoolETH
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.555.18'
$ws.Range("E50").Value = '  +4.21%  '

# Row 51: TheGraph
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").Value = '''0.188'
$ws.Range("E51").Value = '  +15.95%  '

